# Regenerate the localization-status report: the file
# "ead16d7e-00fc-46f2-9b29-ad6b74ad01c9" moved ahead of
# "03fa8877-e67a-4a89-b3b1-39528f7c7118" in the listing (row 4 <-> row 5)
# and ead16d7e's status flipped back to "In Translation" on every sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A4").Value = "ead16d7e-00fc-46f2-9b29-ad6b74ad01c9.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"
$ws.Range("A5").Value = "03fa8877-e67a-4a89-b3b1-39528f7c7118.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"

# --- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A4").Value = "ead16d7e-00fc-46f2-9b29-ad6b74ad01c9.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "ead16d7e-00fc-46f2-9b29-ad6b74ad01c9.1f9f24153240bb293778a52887a8411d810999ac.zh-cn.xlf"
$ws.Range("D4").Value = "2016-02-22 09:02:39"
$ws.Range("A5").Value = "03fa8877-e67a-4a89-b3b1-39528f7c7118.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "03fa8877-e67a-4a89-b3b1-39528f7c7118.140b40af08b626e597f59abd2ea880197ab547dd.zh-cn.xlf"
$ws.Range("D5").Value = "2016-02-22 09:03:20"

# --- de-de sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A4").Value = "ead16d7e-00fc-46f2-9b29-ad6b74ad01c9.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "ead16d7e-00fc-46f2-9b29-ad6b74ad01c9.1f9f24153240bb293778a52887a8411d810999ac.de-de.xlf"
$ws.Range("D4").Value = "2016-02-22 09:02:51"
$ws.Range("A5").Value = "03fa8877-e67a-4a89-b3b1-39528f7c7118.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "03fa8877-e67a-4a89-b3b1-39528f7c7118.140b40af08b626e597f59abd2ea880197ab547dd.de-de.xlf"
$ws.Range("D5").Value = "2016-02-22 09:03:32"
